# Update "想去人数" (want-to-go count) figures in column F across the
# four sheets of the workbook, as produced by the latest gh-pages data
# generation run (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1094
$ws.Range("F11").Value = 403
$ws.Range("F12").Value = 122
$ws.Range("F13").Value = 77
$ws.Range("F15").Value = 438
$ws.Range("F19").Value = 665
$ws.Range("F20").Value = 2528
$ws.Range("F22").Value = 39
$ws.Range("F28").Value = 95
$ws.Range("F32").Value = 33
$ws.Range("F34").Value = 157
$ws.Range("F36").Value = 233

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 711
$ws.Range("F15").Value = 299
$ws.Range("F16").Value = 299
$ws.Range("F19").Value = 922

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2043
$ws.Range("F6").Value = 2254
$ws.Range("F7").Value = 891
$ws.Range("F10").Value = 1085

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2043
$ws.Range("F4").Value = 2254
$ws.Range("F9").Value = 891
$ws.Range("F10").Value = 1085
$ws.Range("F14").Value = 711
$ws.Range("F15").Value = 1094
$ws.Range("F22").Value = 403
$ws.Range("F23").Value = 122
$ws.Range("F25").Value = 77
$ws.Range("F27").Value = 438
$ws.Range("F30").Value = 665
$ws.Range("F31").Value = 2528
$ws.Range("F37").Value = 95
$ws.Range("F41").Value = 299
$ws.Range("F43").Value = 33
$ws.Range("F49").Value = 157
$ws.Range("F51").Value = 233
